$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'E2' '-0.60%'
Set-TextValue $ws 'G2' '19'
Set-TextValue $ws 'D3' '28.55'
Set-TextValue $ws 'E3' '-4.43%'
Set-TextValue $ws 'G3' '19'
Set-TextValue $ws 'D4' '5.237'
Set-TextValue $ws 'E4' '1.15%'
Set-TextValue $ws 'G4' '19'
Set-TextValue $ws 'D5' '0.05695'
Set-TextValue $ws 'E5' '-0.66%'
Set-TextValue $ws 'G5' '19'
Set-TextValue $ws 'E6' '0.37%'
Set-TextValue $ws 'G6' '19'
Set-TextValue $ws 'D7' '3.196'
Set-TextValue $ws 'E7' '3.17%'
Set-TextValue $ws 'G7' '19'
Set-TextValue $ws 'D8' '0.8509'
Set-TextValue $ws 'E8' '-0.65%'
Set-TextValue $ws 'G8' '19'
Set-TextValue $ws 'D9' '0.8566'
Set-TextValue $ws 'E9' '-1.19%'
Set-TextValue $ws 'G9' '19'
Set-TextValue $ws 'D10' '0.1370'
Set-TextValue $ws 'E10' '0.63%'
Set-TextValue $ws 'G10' '19'
Set-TextValue $ws 'D11' '0.07036'
Set-TextValue $ws 'E11' '-0.37%'
Set-TextValue $ws 'G11' '19'
Set-TextValue $ws 'D12' '0.03137'
Set-TextValue $ws 'E12' '7.19%'
Set-TextValue $ws 'G12' '19'
Set-TextValue $ws 'D13' '0.09206'
Set-TextValue $ws 'E13' '-1.94%'
Set-TextValue $ws 'G13' '19'
Set-TextValue $ws 'D14' '0.001536'
Set-TextValue $ws 'E14' '1.70%'
Set-TextValue $ws 'G14' '19'
Set-TextValue $ws 'D15' '0.0005958'
Set-TextValue $ws 'E15' '-0.50%'
Set-TextValue $ws 'G15' '19'
Set-TextValue $ws 'D16' '0.006042'
Set-TextValue $ws 'E16' '0.91%'
Set-TextValue $ws 'G16' '19'
Set-TextValue $ws 'E17' '0.10%'
Set-TextValue $ws 'G17' '19'
Set-TextValue $ws 'E18' '-4.56%'
Set-TextValue $ws 'G18' '19'
Set-TextValue $ws 'E19' '0.43%'
Set-TextValue $ws 'G19' '19'
Set-TextValue $ws 'D20' '0.03272'
Set-TextValue $ws 'E20' '-3.52%'
Set-TextValue $ws 'G20' '19'
Set-TextValue $ws 'D21' '0.1287'
Set-TextValue $ws 'E21' '-2.22%'
Set-TextValue $ws 'G21' '19'
Set-TextValue $ws 'D22' '3.480'
Set-TextValue $ws 'E22' '0.21%'
Set-TextValue $ws 'G22' '19'
Set-TextValue $ws 'D23' '0.04097'
Set-TextValue $ws 'E23' '-1.62%'
Set-TextValue $ws 'G23' '19'
Set-TextValue $ws 'E24' '-0.11%'
Set-TextValue $ws 'G24' '19'
Set-TextValue $ws 'D25' '0.001218'
Set-TextValue $ws 'E25' '-0.52%'
Set-TextValue $ws 'G25' '19'
Set-TextValue $ws 'D26' '0.004138'
Set-TextValue $ws 'E26' '-17.55%'
Set-TextValue $ws 'G26' '19'
Set-TextValue $ws 'E27' '-0.86%'
Set-TextValue $ws 'G27' '19'
Set-TextValue $ws 'D28' '0.0001449'
Set-TextValue $ws 'G28' '19'
Set-TextValue $ws 'G29' '19'
Set-TextValue $ws 'G30' '19'
Set-TextValue $ws 'G31' '19'
Set-TextValue $ws 'G32' '19'
Set-TextValue $ws 'G33' '19'
Set-TextValue $ws 'G34' '19'
Set-TextValue $ws 'G35' '19'
Set-TextValue $ws 'G36' '19'
Set-TextValue $ws 'G37' '19'
Set-TextValue $ws 'G38' '19'
Set-TextValue $ws 'G39' '19'
Set-TextValue $ws 'D40' '0.03765'
Set-TextValue $ws 'E40' '0.35%'
Set-TextValue $ws 'G40' '19'
Set-TextValue $ws 'D41' '0.1063'
Set-TextValue $ws 'E41' '-0.94%'
Set-TextValue $ws 'G41' '19'
Set-TextValue $ws 'D42' '0.003733'
Set-TextValue $ws 'E42' '-35.18%'
Set-TextValue $ws 'G42' '19'
Set-TextValue $ws 'D43' '0.002489'
Set-TextValue $ws 'E43' '24.46%'
Set-TextValue $ws 'G43' '19'
Set-TextValue $ws 'D44' '0.009334'
Set-TextValue $ws 'E44' '-3.42%'
Set-TextValue $ws 'G44' '19'
Set-TextValue $ws 'D45' '0.00005278'
Set-TextValue $ws 'E45' '1.08%'
Set-TextValue $ws 'G45' '19'
Set-TextValue $ws 'E46' '-0.01%'
Set-TextValue $ws 'G46' '19'
Set-TextValue $ws 'D47' '0.07498'
Set-TextValue $ws 'E47' '15.93%'
Set-TextValue $ws 'G47' '19'
Set-TextValue $ws 'D48' '0.002438'
Set-TextValue $ws 'E48' '-3.26%'
Set-TextValue $ws 'G48' '19'
Set-TextValue $ws 'E49' '-0.01%'
Set-TextValue $ws 'G49' '19'
Set-TextValue $ws 'E50' '-0.01%'
Set-TextValue $ws 'G50' '19'
Set-TextValue $ws 'G51' '19'

Write-Host "Applied updates"